$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Laufend" row (old row 8: "Laufend" / "Wenn Distanz > vorherigeDistanz + X")
# and its associated curly-brace shape.
$ws.Shapes.Item("Geschweifte Klammer rechts 3").Delete()
$ws.Rows.Item(8).Delete()

# Remove the "~ 3 sec." row (old row 13, now row 12 after the previous deletion):
# "~ 3 sec." / "laufend" / "stehend" / "Bild machen sobald bewegt (s.o.)" / "Default"
$ws.Rows.Item(12).Delete()

# The remaining "Default" -> "Standard" rename (old G12, now G11).
$ws.Range("G11").Value = "Standard"

# View state: zoom level and selection moved to the "Ablauf" table.
$excel.ActiveWindow.Zoom = 130
[void]$ws.Range("C9:G12").Select()
